$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.769.62"
$ws.Range("E2").Value = "  +0.21%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.535.25"
$ws.Range("E3").Value = "  +0.62%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.64"
$ws.Range("E5").Value = "  -0.18%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.51"
$ws.Range("E6").Value = "  -2.47%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.522"
$ws.Range("E8").Value = "  -1.41%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.533.47"
$ws.Range("E9").Value = "  +0.60%  "

# Row 10
$ws.Range("E10").Value = "  -2.37%  "

# Row 11
$ws.Range("E11").Value = "  +1.01%  "

# Row 12
$ws.Range("E12").Value = "  -0.30%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.02"
$ws.Range("E13").Value = "  -2.80%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.21"
$ws.Range("E14").Value = "  -2.14%  "

# Row 15
$ws.Range("E15").Value = "  +0.19%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000175"
$ws.Range("E16").Value = "  -1.86%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.473.95"
$ws.Range("E17").Value = "  -0.02%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.517.55"
$ws.Range("E18").Value = "  -0.18%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.18"
$ws.Range("E19").Value = "  +2.78%  "

# Row 20
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.76"
$ws.Range("E20").Value = "  +3.01%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "366.01"
$ws.Range("E21").Value = "  +0.98%  "

# Row 22
$ws.Range("B22").Value = "Polkadot"
$ws.Range("C22").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.15"
$ws.Range("E22").Value = "  -0.81%  "

# Row 23
$ws.Range("B23").Value = "Binance-PegBSC-USD"
$ws.Range("C23").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.41"
$ws.Range("E23").Value = "  +42.21%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.52"
$ws.Range("E24").Value = "  -2.49%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.94"
$ws.Range("E25").Value = "  +1.35%  "

# Row 26
$ws.Range("E26").Value = "  +0.02%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.86"
$ws.Range("E27").Value = "  -5.77%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.78"
$ws.Range("E28").Value = "  -4.11%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.662.35"
$ws.Range("E29").Value = "  +0.67%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0944"
$ws.Range("E30").Value = "  -4.59%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "537.14"
$ws.Range("E31").Value = "  -1.91%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.23"
$ws.Range("E32").Value = "  -0.43%  "

# Row 33
$ws.Range("E33").Value = "  -0.35%  "

# Row 34
$ws.Range("E34").Value = "  -4.32%  "

# Row 35
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.128"
$ws.Range("E35").Value = "  -1.91%  "

# Row 36
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.03%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.71"
$ws.Range("E37").Value = "  +3.37%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.30"
$ws.Range("E38").Value = "  +2.86%  "

# Row 39
$ws.Range("E39").Value = "  -2.35%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.63"
$ws.Range("E40").Value = "  +0.07%  "

# Row 41
$ws.Range("E41").Value = "  -1.91%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.06"
$ws.Range("E42").Value = "  -2.33%  "

# Row 43
$ws.Range("E43").Value = "  -3.58%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.17%  "

# Row 45
$ws.Range("E45").Value = "  -3.06%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.37"
$ws.Range("E46").Value = "  -1.24%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.28"
$ws.Range("E47").Value = "  +1.03%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.69"
$ws.Range("E48").Value = "  -0.94%  "

# Row 49
$ws.Range("E49").Value = "  -2.33%  "

# Row 50
$ws.Range("E50").Value = "  -2.76%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.71"
$ws.Range("E51").Value = "  +1.11%  "
